$d = $word.ActiveDocument

$texts = @(
  "במחצית זאת עשינו מלאכות לסוכה, התאמנו על סריגה ועשינו מלאכות יד ותכשיטים לתפארה.`nהייתה אוירה טובה בכיתה וכולן נהנו.`nיהודה אתה ילד נהדר, בהצלחה!",
  "במחצית זאת למדנו את תורת המיספרים, הכרנו כל מספר לעומק, והיתחלנו עם פעולות חשבון בסיסיות,חיבור וחיסור, התקדמנו הרבה עם הספר ""חושבים 1"".והתכוננו לקראת השנה החדשה בההכרה מלמעלה כל כפל וחילוק,`nיהודה אתה ילד נפלא,עלה והצלח!!",
  "במחצית זאת למדנו על עולם המוזיקה, על התווים ועל רמות הקול, התעסקנו עם שירים על מעגל השנה, הייתה אוירה כיפית ונחמדה.`nיהודה אתה תלמיד מדהים!",
  "במחצית זאת למדנו חומש בראשית, למדנו והתפעלנו מבריאת העולם, עקידת יצחק וכו....`nיהודה אתה תלמיד מצוין, בהצלחה!"
)

for ($i = 0; $i -lt $d.Tables.Count; $i++) {
    $table = $d.Tables.Item($i + 1)
    $cell = $table.Cell(1, 2)
    $cell.Range.Text = $texts[$i]
}
